$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 12:23"

$ws.Range("B5").Value = 5996823
$ws.Range("C5").Value = 6242
$ws.Range("E5").Value = 960637
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 94559

$ws.Range("B33").Value = 122673
$ws.Range("C33").Value = 1438
$ws.Range("D33").Value = 98607
$ws.Range("E33").Value = 19348
$ws.Range("G33").Value = 31
$ws.Range("H33").Value = 4718

$ws.Range("B43").Value = 91469
$ws.Range("C43").Value = 851
$ws.Range("D43").Value = 80544
$ws.Range("E43").Value = 10513
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 412

$ws.Range("B44").Value = 90923
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 5880

$ws.Range("B67").Value = 42876
$ws.Range("C67").Value = 662
$ws.Range("D67").Value = 33589
$ws.Range("E67").Value = 8500

$ws.Range("B96").Value = 10919
$ws.Range("C96").Value = 150
$ws.Range("D96").Value = 9835
$ws.Range("E96").Value = 950
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 134

$ws.Range("B97").Value = 10918
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 8749
$ws.Range("E97").Value = 2049
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 120

$ws.Range("B104").Value = 9682
$ws.Range("C104").Value = 105
$ws.Range("D104").Value = 7850
$ws.Range("E104").Value = 1489
$ws.Range("H104").Value = 343

$ws.Range("B105").Value = 9605
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 8385
$ws.Range("E105").Value = 1145
$ws.Range("H105").Value = 75

$ws.Range("B121").Value = 5350
$ws.Range("C121").Value = 159
$ws.Range("D121").Value = 3555
$ws.Range("E121").Value = 1648
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 147

$ws.Range("B122").Value = 5254
$ws.Range("C122").Value = 294
$ws.Range("D122").Value = 1906
$ws.Range("E122").Value = 3320
$ws.Range("H122").Value = 28

$ws.Range("D179").Value = 423
$ws.Range("E179").Value = 37

$ws.Range("B180").Value = 379
$ws.Range("C180").Value = 7
$ws.Range("E180").Value = 38

$ws.Range("B181").Value = 375
$ws.Range("D181").Value = 341
$ws.Range("E181").Value = 34
